$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 50. This pushes the existing rows 50..71
# down to 51..72 (Excel's native Insert shifts all cell data/formatting
# along with the row), matching the diff where the record previously at
# row 50 (and every one after it) ends up one row lower, and the former
# row 71 becomes row 72.
$ws.Rows.Item(50).Insert()

# Populate the newly inserted row 50 with the new weekly record.
$ws.Range("A50").Value = 10
$ws.Range("B50").Value = "Vega Modelo de Temuco"
$ws.Range("C50").Value = "La Araucanía"
$ws.Range("D50").Value = 45016
$ws.Range("E50").Value = 9
$ws.Range("F50").Value = 100112042
$ws.Range("G50").Value = "Locoto"
$ws.Range("H50").Value = "Sin especificar"
$ws.Range("I50").Value = "Primera"
$ws.Range("J50").Value = 80
$ws.Range("K50").Value = 4400
$ws.Range("L50").Value = 4400
$ws.Range("M50").Value = 4400
$ws.Range("N50").Value = "$/kilo"
$ws.Range("O50").Value = "Región de Arica y Parinacota"
$ws.Range("P50").Value = 4400
$ws.Range("Q50").Value = 1
$ws.Range("R50").Value = "Hortaliza"
